# Adds a new "TEMPERATURAS" sheet (as the first tab) with C/F/K conversion
# factors, mirroring the look & feel of the existing "VOLUMENES" sheet, and
# turns the per-cell formulas on VOLUMENES!C2:C13 into a single shared
# formula (as Excel does when the same formula is (re)entered across a
# contiguous range).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet by copying VOLUMENES (so it inherits the exact
#    column widths / header style / page setup) and moving the copy to
#    the very front of the workbook.
# ---------------------------------------------------------------------
$vol = $wb.Worksheets.Item("VOLUMENES")
$firstSheet = $wb.Worksheets.Item(1)
$vol.Copy($firstSheet)

$temps = $wb.Worksheets.Item(1)
$temps.Name = "TEMPERATURAS"

# Drop the extra VOLUMENES data rows (only need header + 6 conversions).
$temps.Rows("8:13").Delete()

# ---------------------------------------------------------------------
# 2. Fill in the temperature conversion symbols / factors / codes.
#    Typed in this order so the shared-string table gets the same
#    insertion order as the source workbook, then sorted alphabetically
#    below (matching the sheet's sortState).
# ---------------------------------------------------------------------
$temps.Range("A2").Value = "C->F"
$temps.Range("A3").Value = "F->C"
$temps.Range("A4").Value = "C->K"
$temps.Range("A5").Value = "F->K"
$temps.Range("A6").Value = "K->C"
$temps.Range("A7").Value = "K->F"

$temps.Range("B2:B7").Value = 1

$temps.Range("C2").Formula = '="rates.put("""&A2&""",new BigDecimal("""&B2&"""));"'
$temps.Range("C3").Formula = '="rates.put("""&A3&""",new BigDecimal("""&B3&"""));"'
$temps.Range("C4").Formula = '="rates.put("""&A4&""",new BigDecimal("""&B4&"""));"'
$temps.Range("C5").Formula = '="rates.put("""&A5&""",new BigDecimal("""&B5&"""));"'
$temps.Range("C6").Formula = '="rates.put("""&A6&""",new BigDecimal("""&B6&"""));"'
$temps.Range("C7").Formula = '="rates.put("""&A7&""",new BigDecimal("""&B7&"""));"'

# Sort A2:C7 by column A (ascending) - same sortState as the other sheets.
$sortRange = $temps.Range("A2:C7")
$temps.Sort.SortFields.Clear()
$temps.Sort.SortFields.Add($temps.Range("A2:A7"))
$temps.Sort.SetRange($sortRange)
$temps.Sort.Header = 0
$temps.Sort.Apply()

$temps.Range("C14").Select()

# ---------------------------------------------------------------------
# 3. VOLUMENES: re-enter the C2:C13 formula as one range write so it is
#    stored as a single shared formula instead of 12 independent ones.
# ---------------------------------------------------------------------
$volumenes = $wb.Worksheets.Item("VOLUMENES")
$volumenes.Range("C2:C13").Formula = '="rates.put("""&A2&""",new BigDecimal("""&B2&"""));"'
